$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-06-14 Saturday"; new="2025-06-15 Sunday"},
    @{old="26×67="; new="43×27="},
    @{old="26×69="; new="56×49="},
    @{old="27×53="; new="66×55="},
    @{old="44×17="; new="20×81="},
    @{old="66×73="; new="37×14="},
    @{old="54×85="; new="65×87="},
    @{old="31×14="; new="69×18="},
    @{old="30×17="; new="75×64="},
    @{old="82×87="; new="95×97="},
    @{old="24×48="; new="53×42="},
    @{old="60×19="; new="95×60="},
    @{old="32×82="; new="12×85="},
    @{old="40×80="; new="58×23="},
    @{old="36×28="; new="86×97="},
    @{old="53×27="; new="14×37="},
    @{old="54×88="; new="16×67="},
    @{old="18×44="; new="92×52="},
    @{old="27×94="; new="59×40="},
    @{old="42×77="; new="20×49="},
    @{old="56×37="; new="22×54="},
    @{old="35×48="; new="22×58="},
    @{old="25×57="; new="64×19="},
    @{old="25×43="; new="26×24="},
    @{old="74×96="; new="26×52="},
    @{old="82×65="; new="46×70="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
